$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12: add B12 part number, update D12 price (formula F12 recalculates automatically)
$ws.Range("B12").Value = "175-5215"
$ws.Range("D12").Value = 77.68

# Row 22: remove G22 value
$ws.Range("G22").ClearContents()

# Row 24: add label for total amazon + shipping
$ws.Range("G24").Value = "Total amazon + frakt"

# Row 31: add label for grand total
$ws.Range("G31").Value = "Sum total kostnad"

# Update selection/view state
$ws.Range("F13").Select()
